$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new values look like plain numbers so Excel
# keeps them as text (matching the original inlineStr cell type).
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D16", "D17", "D19", "D20", "D22", "D23", "D26", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D40", "D41", "D42", "D45", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "42.817.18"
$ws.Range("E2").Value = "  +0.49%  "

# Row 3
$ws.Range("D3").Value = "2.528.71"
$ws.Range("E3").Value = "  -0.41%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").Value = "317.96"
$ws.Range("E5").Value = "  +1.35%  "

# Row 6
$ws.Range("D6").Value = "96.68"
$ws.Range("E6").Value = "  +1.75%  "

# Row 7
$ws.Range("D7").Value = "0.574"
$ws.Range("E7").Value = "  -0.87%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  -0.07%  "

# Row 10
$ws.Range("D10").Value = "35.82"
$ws.Range("E10").Value = "  -0.98%  "

# Row 11
$ws.Range("D11").Value = "0.0816"
$ws.Range("E11").Value = "  +0.09%  "

# Row 12
$ws.Range("D12").Value = "7.53"
$ws.Range("E12").Value = "  -2.17%  "

# Row 13
$ws.Range("E13").Value = "  -4.11%  "

# Row 14
$ws.Range("D14").Value = "2.913.90"
$ws.Range("E14").Value = "  -0.67%  "

# Row 15
$ws.Range("D15").Value = "2.528.68"
$ws.Range("E15").Value = "  -0.17%  "

# Row 16
$ws.Range("D16").Value = "15.04"
$ws.Range("E16").Value = "  -4.08%  "

# Row 17
$ws.Range("D17").Value = "0.848"
$ws.Range("E17").Value = "  -2.11%  "

# Row 18
$ws.Range("D18").Value = "42.834.79"
$ws.Range("E18").Value = "  +0.41%  "

# Row 19
$ws.Range("D19").Value = "6.85"
$ws.Range("E19").Value = "  +2.65%  "

# Row 20
$ws.Range("D20").Value = "12.55"
$ws.Range("E20").Value = "  -3.98%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0964"
$ws.Range("E21").Value = "  -0.60%  "

# Row 22
$ws.Range("D22").Value = "69.62"
$ws.Range("E22").Value = "  -2.05%  "

# Row 23
$ws.Range("D23").Value = "252.61"
$ws.Range("E23").Value = "  -0.88%  "

# Row 24
$ws.Range("E24").Value = "  +0.35%  "

# Row 25
$ws.Range("E25").Value = "  +1.05%  "

# Row 26
$ws.Range("D26").Value = "26.39"
$ws.Range("E26").Value = "  -4.62%  "

# Row 27
$ws.Range("E27").Value = "  -0.04%  "

# Row 28
$ws.Range("D28").Value = "2.41"
$ws.Range("E28").Value = "  +1.84%  "

# Row 29
$ws.Range("D29").Value = "41.04"
$ws.Range("E29").Value = "  +3.41%  "

# Row 30
$ws.Range("E30").Value = "  +3.37%  "

# Row 31
$ws.Range("D31").Value = "5.87"
$ws.Range("E31").Value = "  -1.68%  "

# Row 32
$ws.Range("D32").Value = "156.63"
$ws.Range("E32").Value = "  +0.77%  "

# Row 33
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "2.15"
$ws.Range("E33").Value = "  +0.43%  "

# Row 34
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").Value = "19.42"
$ws.Range("E34").Value = "  -0.54%  "

# Row 35
$ws.Range("D35").Value = "3.38"
$ws.Range("E35").Value = "  -1.07%  "

# Row 36
$ws.Range("E36").Value = "  +3.33%  "

# Row 37
$ws.Range("D37").Value = "0.0790"
$ws.Range("E37").Value = "  -0.18%  "

# Row 38
$ws.Range("D38").Value = "0.113"
$ws.Range("E38").Value = "  +1.28%  "

# Row 39
$ws.Range("E39").Value = "  +9.98%  "

# Row 40
$ws.Range("D40").Value = "0.119"
$ws.Range("E40").Value = "  -0.54%  "

# Row 41
$ws.Range("D41").Value = "21.83"
$ws.Range("E41").Value = "  -11.21%  "

# Row 42
$ws.Range("D42").Value = "0.0305"
$ws.Range("E42").Value = "  +0.78%  "

# Row 43
$ws.Range("E43").Value = "  -1.29%  "

# Row 44
$ws.Range("E44").Value = "  +0.15%  "

# Row 45
$ws.Range("D45").Value = "3.28"
$ws.Range("E45").Value = "  -3.09%  "

# Row 46
$ws.Range("D46").Value = "1.997.82"
$ws.Range("E46").Value = "  -2.64%  "

# Row 47
$ws.Range("D47").Value = "9.13"
$ws.Range("E47").Value = "  +2.27%  "

# Row 48
$ws.Range("D48").Value = "84.17"
$ws.Range("E48").Value = "  -1.14%  "

# Row 49
$ws.Range("D49").Value = "105.79"
$ws.Range("E49").Value = "  +3.54%  "

# Row 50
$ws.Range("D50").Value = "74.98"
$ws.Range("E50").Value = "  -0.54%  "

# Row 51
$ws.Range("D51").Value = "2.770.15"
$ws.Range("E51").Value = "  -0.64%  "

# Remove the temporary text formatting so cell styling matches the original
# (no explicit style index) while keeping the assigned text values intact.
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}